# Fixed the case where we get solution to Dirac equation in 0 iterations.
# Now we can solve O16, Ca40, and Pb208.

$wb = $excel.ActiveWorkbook

# --- Sheet1: N_STEPS 150 -> 300 ---
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("G2").Value = 300

# --- Sheet2: was the selected tab; selection stays G2 ---
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Range("G2").Select()

# --- Sheet3 (spin-up, F column = 0.5): energy guess + 2J+1 values ---
$ws3 = $wb.Worksheets.Item("Sheet3")
$ws3.Range("C2").Formula = "=877"
$ws3.Range("E2").Value = 2
$ws3.Range("E3").Value = 3
$ws3.Range("E4").Value = 3
$ws3.Range("E5").Value = 3
$ws3.Range("D17").Select()

# --- Sheet4 (spin-down, F column = -0.5): energy guess + 2J+1 values ---
$ws4 = $wb.Worksheets.Item("Sheet4")
$ws4.Range("C2").Formula = "=877"
$ws4.Range("E2").Value = 2
$ws4.Range("E3").Value = 3
$ws4.Range("E4").Value = 3
$ws4.Range("E5").Value = 3

# Sheet4 becomes the active / tab-selected sheet, with selection K22
$ws4.Activate()
$ws4.Range("K22").Select()
